# Weekly price update: insert a new "Apio" (celery) price record for
# Terminal Hortofrutícola Agro Chillán as the new row 216, pushing the
# existing rows 216:260 down to 217:261.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 216:260 down to 217:261, leaving a blank row 216 that
# inherits the formatting (incl. the date number format on column D)
# of the row that used to occupy that slot.
$ws.Rows("216:216").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A216").Value = 7
$ws.Range("B216").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C216").Value = "Ñuble"
$ws.Range("D216").Value = 44889
$ws.Range("E216").Value = 16
$ws.Range("F216").Value = 100112017
$ws.Range("G216").Value = "Apio"
$ws.Range("H216").Value = "Americana (o)"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 100
$ws.Range("K216").Value = 9000
$ws.Range("L216").Value = 9500
$ws.Range("M216").Value = 9250
$ws.Range("N216").Value = "`$/docena de matas"
$ws.Range("O216").Value = "Provincia del Elquí"
$ws.Range("P216").Value = 1542
$ws.Range("Q216").Value = 6
$ws.Range("R216").Value = "Hortaliza"
